$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A361:V361").Copy()
$ws.Range("A362:V362").PasteSpecial(-4122)
$ws.Range("A362").Value = "Entrainement"
$ws.Range("B362").Value = 45896
$ws.Range("C362").Value = "Global"
$ws.Range("E362").Value = "Romain Thunet"
$ws.Range("F362").Value = "center back"
$ws.Range("G362").Value = "01:46:05"
$ws.Range("H362").Value = 6.68
$ws.Range("I362").Value = 0.53
$ws.Range("J362").Value = 6.14
$ws.Range("K362").Value = 0.43
$ws.Range("L362").Value = 0.11
$ws.Range("M362").Value = 0
$ws.Range("N362").Value = 0
$ws.Range("O362").Value = 0
$ws.Range("P362").Value = 3.69
$ws.Range("Q362").Value = 24.39
$ws.Range("R362").Value = 4.31
$ws.Range("S362").Value = 32
$ws.Range("T362").Value = 6
$ws.Range("U362").Value = 14
$ws.Range("V362").Value = 5

$ws.Range("A362:V362").Copy()
$ws.Range("A363:V363").PasteSpecial(-4122)
$ws.Range("A363").Value = "Entrainement"
$ws.Range("B363").Value = 45896
$ws.Range("C363").Value = "Global"
$ws.Range("E363").Value = "Levy Ndoutoume"
$ws.Range("F363").Value = "left back"
$ws.Range("G363").Value = "01:46:05"
$ws.Range("H363").Value = 6.47
$ws.Range("I363").Value = 0.48
$ws.Range("J363").Value = 5.97
$ws.Range("K363").Value = 0.38
$ws.Range("L363").Value = 0.09
$ws.Range("M363").Value = 0.02
$ws.Range("N363").Value = 0
$ws.Range("O363").Value = 4
$ws.Range("P363").Value = 3.11
$ws.Range("Q363").Value = 27.31
$ws.Range("R363").Value = 5.46
$ws.Range("S363").Value = 45
$ws.Range("T363").Value = 18
$ws.Range("U363").Value = 37
$ws.Range("V363").Value = 18

$ws.Range("A363:V363").Copy()
$ws.Range("A364:V364").PasteSpecial(-4122)
$ws.Range("A364").Value = "Entrainement"
$ws.Range("B364").Value = 45896
$ws.Range("C364").Value = "Global"
$ws.Range("E364").Value = "Emmanuel Valey"
$ws.Range("F364").Value = "left forward"
$ws.Range("G364").Value = "01:46:21"
$ws.Range("H364").Value = 6.52
$ws.Range("I364").Value = 0.62
$ws.Range("J364").Value = 5.89
$ws.Range("K364").Value = 0.39
$ws.Range("L364").Value = 0.18
$ws.Range("M364").Value = 0.06
$ws.Range("N364").Value = 0
$ws.Range("O364").Value = 6
$ws.Range("P364").Value = 3.6
$ws.Range("Q364").Value = 29.67
$ws.Range("R364").Value = 4.45
$ws.Range("S364").Value = 32
$ws.Range("T364").Value = 5
$ws.Range("U364").Value = 19
$ws.Range("V364").Value = 4

$ws.Range("A364:V364").Copy()
$ws.Range("A365:V365").PasteSpecial(-4122)
$ws.Range("A365").Value = "Entrainement"
$ws.Range("B365").Value = 45896
$ws.Range("C365").Value = "Global"
$ws.Range("E365").Value = "Ilyes Boughanmi"
$ws.Range("F365").Value = "center forward"
$ws.Range("G365").Value = "01:45:49"
$ws.Range("H365").Value = 5.97
$ws.Range("I365").Value = 0.45
$ws.Range("J365").Value = 5.51
$ws.Range("K365").Value = 0.29
$ws.Range("L365").Value = 0.13
$ws.Range("M365").Value = 0.04
$ws.Range("N365").Value = 0
$ws.Range("O365").Value = 4
$ws.Range("P365").Value = 3.31
$ws.Range("Q365").Value = 27.47
$ws.Range("R365").Value = 4.31
$ws.Range("S365").Value = 28
$ws.Range("T365").Value = 1
$ws.Range("U365").Value = 12
$ws.Range("V365").Value = 1

$ws.Range("A365:V365").Copy()
$ws.Range("A366:V366").PasteSpecial(-4122)
$ws.Range("A366").Value = "Entrainement"
$ws.Range("B366").Value = 45896
$ws.Range("C366").Value = "Global"
$ws.Range("E366").Value = "Maé Clavel"
$ws.Range("F366").Value = "left back"
$ws.Range("G366").Value = "01:45:00"
$ws.Range("H366").Value = 6.59
$ws.Range("I366").Value = 0.58
$ws.Range("J366").Value = 6
$ws.Range("K366").Value = 0.46
$ws.Range("L366").Value = 0.12
$ws.Range("M366").Value = 0.01
$ws.Range("N366").Value = 0
$ws.Range("O366").Value = 1
$ws.Range("P366").Value = 3.7
$ws.Range("Q366").Value = 26.38
$ws.Range("R366").Value = 4.67
$ws.Range("S366").Value = 29
$ws.Range("T366").Value = 5
$ws.Range("U366").Value = 11
$ws.Range("V366").Value = 6

$ws.Range("A366:V366").Copy()
$ws.Range("A367:V367").PasteSpecial(-4122)
$ws.Range("A367").Value = "Entrainement"
$ws.Range("B367").Value = 45896
$ws.Range("C367").Value = "Global"
$ws.Range("E367").Value = "Jeremie Laurent"
$ws.Range("F367").Value = "left forward"
$ws.Range("G367").Value = "01:45:01"
$ws.Range("H367").Value = 6.68
$ws.Range("I367").Value = 0.65
$ws.Range("J367").Value = 6.02
$ws.Range("K367").Value = 0.47
$ws.Range("L367").Value = 0.12
$ws.Range("M367").Value = 0.07
$ws.Range("N367").Value = 0
$ws.Range("O367").Value = 5
$ws.Range("P367").Value = 3.76
$ws.Range("Q367").Value = 29.38
$ws.Range("R367").Value = 4.36
$ws.Range("S367").Value = 43
$ws.Range("T367").Value = 2
$ws.Range("U367").Value = 27
$ws.Range("V367").Value = 5

$ws.Range("A367:V367").Copy()
$ws.Range("A368:V368").PasteSpecial(-4122)
$ws.Range("A368").Value = "Entrainement"
$ws.Range("B368").Value = 45896
$ws.Range("C368").Value = "Global"
$ws.Range("E368").Value = "Amine Taiar"
$ws.Range("F368").Value = "center back"
$ws.Range("G368").Value = "01:46:13"
$ws.Range("H368").Value = 6.12
$ws.Range("I368").Value = 0.34
$ws.Range("J368").Value = 5.77
$ws.Range("K368").Value = 0.29
$ws.Range("L368").Value = 0.06
$ws.Range("M368").Value = 0.01
$ws.Range("N368").Value = 0
$ws.Range("O368").Value = 1
$ws.Range("P368").Value = 3.39
$ws.Range("Q368").Value = 26.95
$ws.Range("R368").Value = 4.46
$ws.Range("S368").Value = 26
$ws.Range("T368").Value = 4
$ws.Range("U368").Value = 13
$ws.Range("V368").Value = 4

$ws.Range("A368:V368").Copy()
$ws.Range("A369:V369").PasteSpecial(-4122)
$ws.Range("A369").Value = "Entrainement"
$ws.Range("B369").Value = 45896
$ws.Range("C369").Value = "Global"
$ws.Range("E369").Value = "Omar Benyounes"
$ws.Range("F369").Value = "center midfield"
$ws.Range("G369").Value = "01:45:01"
$ws.Range("H369").Value = 7.38
$ws.Range("I369").Value = 0.81
$ws.Range("J369").Value = 6.55
$ws.Range("K369").Value = 0.65
$ws.Range("L369").Value = 0.15
$ws.Range("M369").Value = 0.03
$ws.Range("N369").Value = 0
$ws.Range("O369").Value = 2
$ws.Range("P369").Value = 4.15
$ws.Range("Q369").Value = 27.53
$ws.Range("R369").Value = 4.26
$ws.Range("S369").Value = 41
$ws.Range("T369").Value = 4
$ws.Range("U369").Value = 27
$ws.Range("V369").Value = 11

$ws.Range("A369:V369").Copy()
$ws.Range("A370:V370").PasteSpecial(-4122)
$ws.Range("A370").Value = "Entrainement"
$ws.Range("B370").Value = 45896
$ws.Range("C370").Value = "Global"
$ws.Range("E370").Value = "Naim Ighbane"
$ws.Range("F370").Value = "center back"
$ws.Range("G370").Value = "01:44:05"
$ws.Range("H370").Value = 5.82
$ws.Range("I370").Value = 0.26
$ws.Range("J370").Value = 5.55
$ws.Range("K370").Value = 0.22
$ws.Range("L370").Value = 0.03
$ws.Range("M370").Value = 0.01
$ws.Range("N370").Value = 0
$ws.Range("O370").Value = 1
$ws.Range("P370").Value = 3.25
$ws.Range("Q370").Value = 26.48
$ws.Range("R370").Value = 4.27
$ws.Range("S370").Value = 26
$ws.Range("T370").Value = 2
$ws.Range("U370").Value = 8
$ws.Range("V370").Value = 4

$ws.Range("A370:V370").Copy()
$ws.Range("A371:V371").PasteSpecial(-4122)
$ws.Range("A371").Value = "Entrainement"
$ws.Range("B371").Value = 45896
$ws.Range("C371").Value = "Global"
$ws.Range("E371").Value = "Karahali Souaré"
$ws.Range("F371").Value = "right forward"
$ws.Range("G371").Value = "01:45:17"
$ws.Range("H371").Value = 7.2
$ws.Range("I371").Value = 0.64
$ws.Range("J371").Value = 6.54
$ws.Range("K371").Value = 0.47
$ws.Range("L371").Value = 0.17
$ws.Range("M371").Value = 0.01
$ws.Range("N371").Value = 0
$ws.Range("O371").Value = 4
$ws.Range("P371").Value = 3.69
$ws.Range("Q371").Value = 26.64
$ws.Range("R371").Value = 5.21
$ws.Range("S371").Value = 54
$ws.Range("T371").Value = 10
$ws.Range("U371").Value = 33
$ws.Range("V371").Value = 10

$ws.Range("A371:V371").Copy()
$ws.Range("A372:V372").PasteSpecial(-4122)
$ws.Range("A372").Value = "Entrainement"
$ws.Range("B372").Value = 45896
$ws.Range("C372").Value = "Global"
$ws.Range("E372").Value = "Mattheo Haon"
$ws.Range("F372").Value = "right back"
$ws.Range("G372").Value = "01:45:41"
$ws.Range("H372").Value = 7.23
$ws.Range("I372").Value = 0.69
$ws.Range("J372").Value = 6.53
$ws.Range("K372").Value = 0.47
$ws.Range("L372").Value = 0.21
$ws.Range("M372").Value = 0.02
$ws.Range("N372").Value = 0
$ws.Range("O372").Value = 4
$ws.Range("P372").Value = 4.05
$ws.Range("Q372").Value = 27.55
$ws.Range("R372").Value = 4.29
$ws.Range("S372").Value = 29
$ws.Range("T372").Value = 5
$ws.Range("U372").Value = 14
$ws.Range("V372").Value = 10

$ws.Range("A372:V372").Copy()
$ws.Range("A373:V373").PasteSpecial(-4122)
$ws.Range("A373").Value = "Entrainement"
$ws.Range("B373").Value = 45896
$ws.Range("C373").Value = "Global"
$ws.Range("E373").Value = "Yoan Zouma"
$ws.Range("F373").Value = "center back"
$ws.Range("G373").Value = "01:44:52"
$ws.Range("H373").Value = 5.7
$ws.Range("I373").Value = 0.29
$ws.Range("J373").Value = 5.4
$ws.Range("K373").Value = 0.22
$ws.Range("L373").Value = 0.07
$ws.Range("M373").Value = 0.01
$ws.Range("N373").Value = 0
$ws.Range("O373").Value = 1
$ws.Range("P373").Value = 3.09
$ws.Range("Q373").Value = 27.28
$ws.Range("R373").Value = 4.58
$ws.Range("S373").Value = 12
$ws.Range("T373").Value = 2
$ws.Range("U373").Value = 9
$ws.Range("V373").Value = 5

$ws.Range("A373:V373").Copy()
$ws.Range("A374:V374").PasteSpecial(-4122)
$ws.Range("A374").Value = "Entrainement"
$ws.Range("B374").Value = 45896
$ws.Range("C374").Value = "Global"
$ws.Range("E374").Value = "Ilan Ihaddadene"
$ws.Range("F374").Value = "center midfield"
$ws.Range("G374").Value = "01:44:45"
$ws.Range("H374").Value = 7.52
$ws.Range("I374").Value = 0.42
$ws.Range("J374").Value = 7.1
$ws.Range("K374").Value = 0.39
$ws.Range("L374").Value = 0.03
$ws.Range("M374").Value = 0.01
$ws.Range("N374").Value = 0
$ws.Range("O374").Value = 1
$ws.Range("P374").Value = 4.25
$ws.Range("Q374").Value = 25.92
$ws.Range("R374").Value = 5.21
$ws.Range("S374").Value = 35
$ws.Range("T374").Value = 5
$ws.Range("U374").Value = 14
$ws.Range("V374").Value = 2

$ws.Range("D362:D374").Value = "J-3"

$ws.Application.CutCopyMode = $false
$ws.Range("E377").Select()
